# "Timeslot read error fixed"
# The extra demo accounts (admin2..admin6) and their matching timeslot rows
# are removed, the remaining "admin" account gets a regenerated ID/DOB, and
# the Timeslot Information sheet's data row is updated to match that ID.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: Account Information ---------------------------------------
$ws1 = $wb.Worksheets.Item("Account Information")

# Drop the admin2..admin6 rows (rows 3-7), keeping only the header + admin.
$ws1.Rows("3:7").Delete()

$ws1.Range("A2").Value = "admin"
$ws1.Range("B2").Value = "admin"
$ws1.Range("C2").Value = "admin@gmail.com"
$ws1.Range("D2").Value = "Admin"

# Force DateOfBirth to stay plain text ("11/11/1111") instead of being
# auto-parsed into a date serial, then drop the quote-prefix style so the
# cell lands back on the default style like the rest of the sheet.
$ws1.Range("E2").Value = "'11/11/1111"
$ws1.Range("E2").Style = "Normal"

$ws1.Range("F2").Value = 48914717

$ws1.Range("A2:F7").Select()

# --- Sheet 2: Timeslot Information ---------------------------------------
$ws2 = $wb.Worksheets.Item("Timeslot Information")

# Drop the matching extra timeslot rows (rows 3-7).
$ws2.Rows("3:7").Delete()

$ws2.Range("A2").Value = 48914717
$ws2.Range("B2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("C2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("D2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("E2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("F2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("G2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"
$ws2.Range("H2").Value = "null,null,null,null,Biology,Biology,Biology,Biology,Biology,Biology,Biology,Biology"

$ws2.Range("I16:I17").Select()

# --- Sheet 3: Class Information -------------------------------------------
$ws3 = $wb.Worksheets.Item("Class Information")
$ws3.Range("F20").Select()

# Re-activate the Class Information tab, matching activeTab="2" in the book.
$ws3.Activate()
